# Update the Mapping sheet's reserve/station bounding-box coordinates after
# converting all shapefiles to WGS 84 (EPSG 4269).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# New bounding-box longitude/latitude pairs (Res_Bounding_Box / SK_Bounding_Box)
$ws.Range("A2").Value = -70.5611
$ws.Range("B2").Value = -70.5545
$ws.Range("A3").Value = 41.5361
$ws.Range("B3").Value = 41.541
$ws.Range("A4").Value = -70.4697
$ws.Range("B4").Value = -70.4764
$ws.Range("A5").Value = 41.65
$ws.Range("B5").Value = 41.6451

# Touch previously-unused cells so they become present (but empty) entries,
# matching the wider rectangular extent used by the refreshed data paste.
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("A6:F6").Style = "Normal"
